$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 30.31806102698224
$ws.Range("B3").Value = 30.76153967020114
$ws.Range("B4").Value = 40.25973713831916
$ws.Range("H5").Value = 96.24236789370403
$ws.Range("H6").Value = 96.13098304617974
$ws.Range("H7").Value = 96.22756446181195
$ws.Range("C8").Value = 50.63498013815493
$ws.Range("C9").Value = 50.43470392334577
$ws.Range("C10").Value = 50.70436496210364
$ws.Range("D11").Value = 98.69172919731389
$ws.Range("D12").Value = 98.78537891235521
$ws.Range("D13").Value = 98.97802536907948
$ws.Range("E14").Value = 98.93984262359832
$ws.Range("E15").Value = 98.9978291233006
$ws.Range("E16").Value = 98.93757540209221
$ws.Range("F17").Value = 98.4813239258053
$ws.Range("F18").Value = 98.58154259493311
$ws.Range("F19").Value = 98.44626910233967
$ws.Range("G20").Value = 97.56959798991242
$ws.Range("G21").Value = 97.61725417075434
$ws.Range("G22").Value = 97.61390091135038
$ws.Range("B23").Value = 40.98592402049587
$ws.Range("B24").Value = 49.03583991382323
$ws.Range("H25").Value = 96.30820249441508
$ws.Range("H26").Value = 96.23109626382599
$ws.Range("C27").Value = 50.19379561441237
$ws.Range("C28").Value = 50.30764186347617
$ws.Range("D29").Value = 98.8470680953949
$ws.Range("D30").Value = 98.79434438846096
$ws.Range("E31").Value = 98.95227071897385
$ws.Range("E32").Value = 98.95128926113065
$ws.Range("F33").Value = 98.45407170166315
$ws.Range("F34").Value = 98.58918864167788
$ws.Range("G35").Value = 97.50497962783544
$ws.Range("G36").Value = 97.62698278886836
$ws.Range("B37").Value = 36.33344537134072
$ws.Range("B38").Value = 50.80760093129585
$ws.Range("H39").Value = 96.28893688435859
$ws.Range("H40").Value = 96.09449641807879
$ws.Range("C41").Value = 50.94024582960805
$ws.Range("C42").Value = 50.10465028580613
$ws.Range("D43").Value = 98.9599121273158
$ws.Range("D44").Value = 98.79257303335514
$ws.Range("E45").Value = 98.98878201414485
$ws.Range("E46").Value = 98.85244615087099
$ws.Range("F47").Value = 98.57621715884011
$ws.Range("F48").Value = 98.44329350150969
$ws.Range("G49").Value = 97.49120598448928
$ws.Range("G50").Value = 97.63882725993936
